$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Markdown text for "question1 daily" -> replaced with markdown question 1 body
$mdQ1 = '#MarkDown Question 1' + $nl + '* point1' + $nl + '``` Code code ```' + $nl + 'Answer - abc'

# Markdown text for "question2 daily" -> replaced with markdown question 2 body
$mdQ2 = '#MarkDown Question 2' + $nl + '``` Code code ```' + $nl + 'Answer - bcd'

# Update existing row 2 (was "question1 daily")
$ws.Range("B2").Value = $mdQ1
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 57.6

# Update existing row 4 (was "question2 daily")
$ws.Range("B4").Value = $mdQ2
$ws.Range("B4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 43.2

# Add new row 7, duplicate of the (updated) row 2 content
$ws.Range("A7").Value = $true
$ws.Range("B7").Value = $mdQ1
$ws.Range("C7").Value = "abc"
$ws.Range("D7").Value = $true
$ws.Range("B7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 57.6

# Add new row 8, duplicate of the (updated) row 4 content
$ws.Range("A8").Value = $true
$ws.Range("B8").Value = $mdQ2
$ws.Range("C8").Value = "bcd"
$ws.Range("D8").Value = $true
$ws.Range("B8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 43.2

# Match the final selection state from the diff
$ws.Range("A7").Select()
